$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the "Marking" correct-answer score and recompute the total/max display
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 115
$ws.Range("E12").Value = "115/140"
